# Auto-update draw results: append the 2025-11-14 Pick 4 draw as a new
# row at the bottom of the "Results" sheet (row 59), mirroring the daily
# automated export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 59

# Columns A (date) and C (phase code) look numeric/date-like to Excel's
# auto-detection, so force them to text first so the literal strings are
# preserved exactly as produced by the upstream export (e.g. "2025-11-14"
# and "251114" rather than being coerced into date/number values).
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("C" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-11-14"
$ws.Range("B" + $newRow).Value = "Pick 4"
$ws.Range("C" + $newRow).Value = "251114"
$ws.Range("D" + $newRow).Value = "8-0-5-9"
$ws.Range("E" + $newRow).Value = "2025-11-14T21:39:22.299+04:00"
